# Rename + reposition/resize the dashboard background shapes on slide 1
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$sh = $s.Shapes.Item(1)  # id=33
$sh.Name = "Slicer Panel Background"
$sh.Left = 78.49236220472442
$sh.Top = 44.07787513917638
$sh.Width = 757.6929133858267
$sh.Height = 42.42653543307087

$sh = $s.Shapes.Item(2)  # id=28
$sh.Name = "Title Bar"
$sh.Left = 19.35779527559055
$sh.Top = 7.0346457958916035
$sh.Width = 806.1594543730969
$sh.Height = 28.80000019168854

$sh = $s.Shapes.Item(3)  # id=29
$sh.Name = "Dashboard Title Left Accent"
$sh.Left = 8.689685345022697
$sh.Top = 7.0346457958916035
$sh.Width = 10.668110370688552
$sh.Height = 28.80000019168854

$sh = $s.Shapes.Item(4)  # id=30
$sh.Name = "Dashboard Title Bar Right Accent"
$sh.Left = 825.5172440944882
$sh.Top = 7.0346457958916035
$sh.Width = 10.668110370688552
$sh.Height = 28.80000019168854

$sh = $s.Shapes.Item(5)  # id=31
$sh.Name = "Navigation Panel Left Side"
$sh.Left = 3.5286614173228346
$sh.Top = 44.07779527559055
$sh.Width = 114.81622047244095
$sh.Height = 417.10331726422646

$sh = $s.Shapes.Item(6)  # id=32
$sh.Name = "Navigation Panel Right Side"
$sh.Left = 76.13795275590552
$sh.Top = 44.07779527559055
$sh.Width = 45.06
$sh.Height = 417.10331726422646

$sh = $s.Shapes.Item(7)  # id=34
$sh.Name = "Total Card"
$sh.Left = 134.13354330708663
$sh.Top = 98.13031496062992
$sh.Width = 180.41378021420456
$sh.Height = 86.4

$sh = $s.Shapes.Item(8)  # id=35
$sh.Name = "Province Card"
$sh.Left = 131.58622047244094
$sh.Top = 192.00008392814576
$sh.Width = 180.41378021420456
$sh.Height = 109.24126052916519

$sh = $s.Shapes.Item(9)  # id=36
$sh.Name = "Postcode Card"
$sh.Left = 131.5861434946736
$sh.Top = 313.6549606299213
$sh.Width = 707.7517395161306
$sh.Height = 147.52614593607984

$sh = $s.Shapes.Item(10)  # id=37
$sh.Name = "District Card"
$sh.Left = 324.9180450458677
$sh.Top = 98.13031496062992
$sh.Width = 252.0
$sh.Height = 206.33850860619734

$sh = $s.Shapes.Item(11)  # id=38
$sh.Name = "Sub District Card"
$sh.Left = 587.2888488964171
$sh.Top = 98.13031496062992
$sh.Width = 252.0
$sh.Height = 206.33850860619734

# Refresh the cached "datetimeFigureOut" placeholder text (master + every layout)
# from 7/26/2023 -> 7/28/2023, matching the date the deck was re-saved.
$newDate = "7/28/2023"
$master = $p.SlideMaster

for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $phsh = $master.Shapes.Item($j)
    if ($phsh.Name -like "Date Placeholder*") {
        $phsh.TextFrame.TextRange.Text = $newDate
    }
}

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $cl = $master.CustomLayouts.Item($i)
    for ($j = 1; $j -le $cl.Shapes.Count; $j++) {
        $phsh = $cl.Shapes.Item($j)
        if ($phsh.Name -like "Date Placeholder*") {
            $phsh.TextFrame.TextRange.Text = $newDate
        }
    }
}
